# Word COM-interop script implementing the edit described by the diff.
# Strategy:
#  - Simple full-paragraph text substitutions are done with Find/Replace
#    (Content.Find.Execute), using long, unique anchor strings so each
#    match is unambiguous.
#  - New paragraphs are inserted with Range.InsertParagraphAfter() on an
#    anchor paragraph located by scanning Paragraphs for unique text,
#    then the new (empty) paragraph is reached via .Next() and populated.
#  - The renamed heading bookmark is recreated (delete + add) so its
#    Name matches the new heading text/slug.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. Title: "La produttività del lavoro in Italia non cresce." ->
#           "Formazione, salari e produttività."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "La produttività del lavoro in Italia non cresce.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Formazione, salari e produttività.", 2)

# ---------------------------------------------------------------------
# 2. Big intro paragraph (FirstParagraph style): rewrite the three
#    sentences that changed, leaving the opening sentence untouched.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Questi fattori messi assieme generano la crescita di produttività del lavoro e la ricchezza da distribuire (inclusa la crescita dei salari).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Questi fattori messi assieme generano la crescita della produttività del lavoro e soprattutto generano la ricchezza da distribuire (anche attraverso crescita dei salari reali).", 2)

$d.Content.Find.Execute(
    "Dove non ci sono investimenti e innovazione di modelli organizzativi la formazione resta confinata fra le attività che servono a provare la buona volontà di un assistito dal welfare statale, nei regimi in cui la percezione di sussidi è sottoposta a condizioni (fra le quali svolgere politiche attive).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dove non ci sono investimenti e innovazione nei modelli organizzativi la formazione resta confinata fra le attività che servono a provare la buona volontà di un assistito dal welfare statale; nei paesi in cui la percezione di sussidi è sottoposta a condizioni (fra le quali svolgere politiche attive) chi non si forma o non si cerca un lavoro non riceve sussidi.", 2)

$d.Content.Find.Execute(
    "Per questo esistono al mondo regimi diversi di formazione e di utilizzo di istituti come l’apprendistato: dipendono dal grado di innovazione e dal tasso di crescita del sistema economico.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Per questo esistono nel mondo diversi regimi del contratto di apprendistato e modalità diverse di organizzazione della formazione: dipendono dal grado di innovazione e dal tasso di crescita del sistema economico.", 2)

# ---------------------------------------------------------------------
# 3. "contabilità nazionali" -> "contabilità nazionale"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "sono una rielaborazione dei principali indicatori di contabilità nazionali, armonizzati",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sono una rielaborazione dei principali indicatori di contabilità nazionale, armonizzati", 2)

# ---------------------------------------------------------------------
# 4. "La figura 1 illustra..." paragraph: rewrite the second half.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Si tratta di misure comparabili e che tengono conto dei diversi livelli di prezzo, ma dato che sono medie non tengono conto di eventuali variazioni nella distribuzione dei redditi nel tempo. Il grafico mostra minimo, massimo e ultimo valore disponibile.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Si tratta di misure comparabili e che tengono conto dei diversi livelli di inflazione, ma dato che sono medie non tengono conto di eventuali variazioni nella distribuzione dei redditi fra i lavoratori all’interno di ogni paese. Il grafico mostra il salario medio annuo minimo nei 25 anni presi in considerazione, il massimo e l’ultimo valore disponibile.", 2)

# ---------------------------------------------------------------------
# 5. New paragraph "Nella figura 1bis..." inserted right after the
#    "I salari italiani sono in discesa..." paragraph (before the
#    salario_medio-1.png picture paragraph).
# ---------------------------------------------------------------------
$anchor1 = Find-ParagraphByText $d "I salari italiani sono in discesa dal 2010"
$anchor1.Range.InsertParagraphAfter()
$newPara1 = $anchor1.Next()
$newPara1.Range.Text = "Nella figura 1bis vediamo le stesse curve in 6 grafici affiancati che hanno la stessa scala. In questo modo possiamo osservare le differenze fra i livelli dei salari medi ed il loro andamento nel tempo."
$newPara1.Style = "BodyText"

# ---------------------------------------------------------------------
# 6. New paragraph "La curva dei salari italiani..." inserted right
#    after the salario_medio-1.png picture paragraph (which now follows
#    $newPara1) and before the "Produttività: il Pil per ora lavorata."
#    heading.
# ---------------------------------------------------------------------
$pictureAfterPara1 = $newPara1.Next()
$pictureAfterPara1.Range.InsertParagraphAfter()
$newPara2 = $pictureAfterPara1.Next()
$newPara2.Range.Text = "La curva dei salari italiani è la più bassa, quella dei salari negli USA la più alta. Italia e Giappone hanno curve sostanzialmente piatte. Annotiamo questa similarità fra Italia e Giappone, la riprenderemo in seguito."
$newPara2.Style = "BodyText"

# ---------------------------------------------------------------------
# 7. "Se invece del PIL per ora lavorata..." paragraph: rewrite ending.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "visto che per tutti i paesi la curva passa a quota 100 nel 2015, ma come dinamica relativa.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "visto che per tutti i paesi la curva passa a quota 100 nel 2015, ma sono utili per osservare la dinamica relativa del PIL per addetto nel tempo per ogni singolo paese.", 2)

# ---------------------------------------------------------------------
# 8. "Questo indicatore di produttività..." paragraph: CIG sentence
#    rewrite.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "la principale risposta in termini di politica del lavoro sia stata la Cassa Integrazione Guadagni (CIG) che mantiene il numero degli occupati, ma fa scendere le ore lavorate. E’ chiaro l’impatto sul rapporto PIL per addetto nel caso di uso intensivo di CIG.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "la principale risposta in termini di politica del lavoro sia stata la Cassa Integrazione Guadagni (CIG) che ha mantenuto relativamente più elevato il numero degli occupati, ma che ha fatto scendere le ore lavorate. E’ chiaro l’impatto negativo sul rapporto PIL per addetto nel caso di uso intensivo di CIG.", 2)

# ---------------------------------------------------------------------
# 9. "Si, è illustrata dalla figura 3..." -> "Si, esiste ed è illustrata..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Si, è illustrata dalla figura 3.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Si, esiste ed è illustrata dalla figura 3.", 2)

# ---------------------------------------------------------------------
# 10. Heading "Alcune conclusioni" -> "Alcune sottolineature conclusive."
#     plus renaming its bookmark.
# ---------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("alcune-conclusioni")
$bookmarkRange = $oldBookmark.Range
$oldBookmark.Delete()

$d.Content.Find.Execute(
    "Alcune conclusioni",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Alcune sottolineature conclusive.", 2)

$headingPara = Find-ParagraphByText $d "Alcune sottolineature conclusive."
$d.Bookmarks.Add("alcune-sottolineature-conclusive.", $headingPara.Range)

# ---------------------------------------------------------------------
# 11. "...mi limito a commentare le politiche del governo giapponese..."
#     -> "...mi limito a una breve sintesi delle politiche del governo
#     giapponese..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Visto che le politiche italiane sono note mi limito a commentare le politiche del governo giapponese:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Visto che le politiche italiane sono note mi limito a una breve sintesi delle politiche del governo giapponese:", 2)
